$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("C8").Value = 7

$ws.Range("D8").Value = "'2"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "'Short point (up to 3 mtr.)"
$ws.Range("E8").Style = "Normal"

$ws.Range("F8").Value = 256

$ws.Range("G8").Value = "'1792.00"
$ws.Range("G8").Style = "Normal"

# Row 9
$ws.Range("C9").Value = 96

$ws.Range("G9").Value = "'63552.00"
$ws.Range("G9").Style = "Normal"

# Row 10
$ws.Range("A10").Value = "'"
$ws.Range("A10").Style = "Normal"

$ws.Range("C10").Value = 75

$ws.Range("D10").Value = "'2.0"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "'Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it's  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet's & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("E10").Style = "Normal"

$ws.Range("F10").Value = 0

$ws.Range("G10").Value = "'0.00"
$ws.Range("G10").Style = "Normal"

# Row 11
$ws.Range("C11").Value = 65

$ws.Range("D11").Value = "'4.0"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("E11").Style = "Normal"

$ws.Range("F11").Value = 50

$ws.Range("G11").Value = "'3250.00"
$ws.Range("G11").Style = "Normal"

# Row 12
$ws.Range("A12").Value = "'"
$ws.Range("A12").Style = "Normal"

$ws.Range("C12").Value = 81

$ws.Range("D12").Value = "'12.0"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "'Supplying and drawing FR PVC insulated & unsheathed flexible copper conductor as per PWD specification for electrical Works with ISI marked (IS:694) and as per IS 8130 : 2013 of 1.1 kV grade . Wire should be made from  99.90 % purity copper, class 2 stranding in acc. to IS:8130/IEC 60228 for  lower watt loss , oxygen free for less chances of oxidization, insulation PVC type A/C/D , flame retardant as per IS 10810-53, better amperage rating as per IS:3961 part 5,  in existing  surface or recessed PVC/ MS conduit/casing capping making connections with Copper Lugs of suitable size, Ferrules,testing etc. as required. OEM Must have its own in house NABL lab setup for all testing facilities for wires.   For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("E12").Style = "Normal"

$ws.Range("F12").Value = 0

$ws.Range("G12").Value = "'0.00"
$ws.Range("G12").Style = "Normal"

# Row 13
$ws.Range("A13").Value = "'Mtr."
$ws.Range("A13").Style = "Normal"

$ws.Range("C13").Value = 71

$ws.Range("D13").Value = "'19"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "'2 x 2.5 sq. mm. + 1x1.5sqmm"
$ws.Range("E13").Style = "Normal"

$ws.Range("F13").Value = 81

$ws.Range("G13").Value = "'5751.00"
$ws.Range("G13").Style = "Normal"

# Row 14
$ws.Range("A14").Value = "'Set"
$ws.Range("A14").Style = "Normal"

$ws.Range("C14").Value = 48

$ws.Range("D14").Value = "'13.0"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. 'B' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR .   "
$ws.Range("E14").Style = "Normal"

$ws.Range("F14").Value = 5733

$ws.Range("G14").Value = "'275184.00"
$ws.Range("G14").Style = "Normal"

# Row 15
$ws.Range("A15").Value = "'Each"
$ws.Range("A15").Style = "Normal"

$ws.Range("C15").Value = 71

$ws.Range("D15").Value = "'25"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "'1200 mm Sweep BEE 1 Star rated (service value >=4.0 to < 4.5 )"
$ws.Range("E15").Style = "Normal"

$ws.Range("F15").Value = 1890

$ws.Range("G15").Value = "'134190.00"
$ws.Range("G15").Style = "Normal"

# Row 16
$ws.Range("C16").Value = 12

$ws.Range("D16").Value = "'16.0"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("A17").Value = "'Each"
$ws.Range("A17").Style = "Normal"

$ws.Range("C17").Value = 87

$ws.Range("D17").Value = "'27"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "'1170mm(+/-10%) LED batten with min. lumen output 2200 lm"
$ws.Range("E17").Style = "Normal"

$ws.Range("F17").Value = 492

$ws.Range("G17").Value = "'42804.00"
$ws.Range("G17").Style = "Normal"

# Row 18
$ws.Range("C18").Value = 29

# Row 19
$ws.Range("A19").Value = "'"
$ws.Range("A19").Style = "Normal"

$ws.Range("C19").Value = 96

$ws.Range("D19").Value = "'31"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "'Double pole MCB(With B/C curve tripping Characteristics)"
$ws.Range("E19").Style = "Normal"

$ws.Range("F19").Value = 0

$ws.Range("G19").Value = "'0.00"
$ws.Range("G19").Style = "Normal"

# Row 20
$ws.Range("A20").Value = "'Each"
$ws.Range("A20").Style = "Normal"

$ws.Range("C20").Value = 89

$ws.Range("D20").Value = "'32"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "' 50/63 A rating"
$ws.Range("E20").Style = "Normal"

$ws.Range("F20").Value = 900

$ws.Range("G20").Value = "'80100.00"
$ws.Range("G20").Style = "Normal"

# Row 21
$ws.Range("A21").Value = "'"
$ws.Range("A21").Style = "Normal"

$ws.Range("C21").Value = 83

$ws.Range("D21").Value = "'18.0"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("E21").Style = "Normal"

$ws.Range("F21").Value = 0

$ws.Range("G21").Value = "'0.00"
$ws.Range("G21").Style = "Normal"

# Row 22
$ws.Range("A22").Value = "'%"
$ws.Range("A22").Style = "Normal"

$ws.Range("C22").Value = 6

$ws.Range("D22").Value = "'37"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "'Add Tender Premium "
$ws.Range("E22").Style = "Normal"

# Row 24
$ws.Range("G24").Value = "'606623.00"
$ws.Range("G24").Style = "Normal"

$ws.Range("H24").Value = "'606623.00"
$ws.Range("H24").Style = "Normal"

# Row 26
$ws.Range("G26").Value = "'606623.00"
$ws.Range("G26").Style = "Normal"

$ws.Range("H26").Value = "'606623.00"
$ws.Range("H26").Style = "Normal"
